# Update countries & provincias Spain
# Daily refresh of the COVID-19 country stats table:
#  - Most rows get refreshed B:H figures (new totals/cases/recoveries/deaths).
#  - Egipto overtakes Tailandia (rows 53/54), and Irak overtakes Estonia
#    (rows 68/69) in the case-count ranking, so those two row pairs swap
#    which country they show while carrying their own refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> values for columns A..H
$updates = @{
    4  = @("Estados Unidos", 685541, 7971, 58160, 591881, 13369, 883, 35500)
    6  = @("Italia",         172434, 3493, 42727, 106962, 2812,  575, 22745)
    15 = @("Brasil",         30961,  278,  14026, 14979,  6634,  9,   1956)
    16 = @("Canada",         30697,  591,  10092, 19353,  557,   57,  1252)
    44 = @("Singapur",       5050,   623,  708,   4331,   22,    1,   11)
    53 = @("Egipto",         2844,   171,  646,   1993,   0,     9,   205)
    54 = @("Tailandia",      2700,   28,   1689,  964,    61,    1,   47)
    68 = @("Irak",           1482,   48,   906,   495,    0,     1,   81)
    69 = @("Estonia",        1459,   25,   145,   1276,   11,    2,   38)
    82 = @("Cuba",           923,    61,   192,   700,    16,    4,   31)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($col = 0; $col -lt $vals.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $vals[$col]
    }
}
